$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the obsolete "1-RAP / Squilla mantis" row (row 10); all data below
# (rows 11-52) shifts up by one, so the sheet's used range shrinks to K51.
$ws.Rows("10").Delete()
